$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.900.73'
$ws.Range("E2").Value = '  +2.04%  '
$ws.Range("D3").Value = '1.811.85'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("D4").Value = "'" + '1.006'
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("D5").Value = "'" + '313.86'
$ws.Range("E5").Value = '  +3.18%  '
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("D8").Value = "'" + '0.3697'
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").Value = "'" + '0.07242'
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("D10").Value = "'" + '0.8661'
$ws.Range("E10").Value = '  +4.12%  '
$ws.Range("D11").Value = '2.104.56'
$ws.Range("E11").Value = '  +18.51%  '
$ws.Range("D12").Value = "'" + '21.31'
$ws.Range("E12").Value = '  +5.78%  '
$ws.Range("D13").Value = "'" + '6.632'
$ws.Range("E13").Value = '  +3.83%  '
$ws.Range("E14").Value = '  +3.43%  '
$ws.Range("D15").Value = "'" + '0.06935'
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("D16").Value = "'" + '80.85'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Value = "'" + '1.006'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = "'" + '0.000008831'
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = "'" + '15.29'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").Value = '26.942.49'
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("D22").Value = "'" + '5.192'
$ws.Range("E22").Value = '  +3.93%  '
$ws.Range("D23").Value = "'" + '10.98'
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").Value = '2.346.01'
$ws.Range("E24").Value = '  +17.88%  '
$ws.Range("D25").Value = "'" + '154.32'
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("D26").Value = "'" + '1.887'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").Value = "'" + '18.38'
$ws.Range("E27").Value = '  +1.73%  '
$ws.Range("D28").Value = "'" + '5.238'
$ws.Range("E28").Value = '  +4.25%  '
$ws.Range("D29").Value = "'" + '1.941'
$ws.Range("E29").Value = '  +16.39%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = "'" + '0.08945'
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("D32").Value = "'" + '1.167'
$ws.Range("E32").Value = '  +4.76%  '
$ws.Range("D33").Value = "'" + '0.7438'
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").Value = "'" + '4.435'
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("E35").Value = '  +3.31%  '
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").Value = "'" + '1.122'
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("E38").Value = '  +3.10%  '
$ws.Range("D39").Value = "'" + '0.01927'
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("D41").Value = "'" + '2.756'
$ws.Range("E41").Value = '  +10.73%  '
$ws.Range("D42").Value = "'" + '0.1650'
$ws.Range("E42").Value = '  +3.14%  '
$ws.Range("D43").Value = "'" + '6.496'
$ws.Range("E43").Value = '  +4.94%  '
$ws.Range("D44").Value = "'" + '8.303'
$ws.Range("E44").Value = '  +4.09%  '
$ws.Range("D45").Value = "'" + '107.31'
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").Value = "'" + '10.40'
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D50").Value = "'" + '0.06283'
$ws.Range("E50").Value = '  +1.62%  '
$ws.Range("D51").Value = "'" + '1.802'
$ws.Range("E51").Value = '  +4.53%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'" + '1.651'
$ws.Range("E48").Value = '  +5.34%  '
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").Value = "'" + '0.4563'
$ws.Range("E49").Value = '  +2.08%  '
